# Auto-generated script to update 2024 (and a couple 2021/2022) crime counts
# to reflect newly added data for 2024-10-24 across the workbook's sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('K2').Value = 6601
$ws.Range('K3').Value = 6821
$ws.Range('H4').Value = 1740
$ws.Range('I4').Value = 1807
$ws.Range('K4').Value = 1414
$ws.Range('K5').Value = 495
$ws.Range('K6').Value = 7496
$ws.Range('H7').Value = 26053
$ws.Range('I7').Value = 26267
$ws.Range('K7').Value = 22827

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('K6').Value = 123
$ws.Range('K7').Value = 289

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('K3').Value = 458
$ws.Range('K4').Value = 86
$ws.Range('K6').Value = 497
$ws.Range('K7').Value = 1497

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('K2').Value = 170
$ws.Range('K3').Value = 176
$ws.Range('K6').Value = 113
$ws.Range('K7').Value = 494

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range('K2').Value = 127
$ws.Range('K3').Value = 132
$ws.Range('K7').Value = 375

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('K2').Value = 221
$ws.Range('K3').Value = 256
$ws.Range('K7').Value = 775

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('K2').Value = 99
$ws.Range('K3').Value = 159
$ws.Range('K6').Value = 96
$ws.Range('K7').Value = 384

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('K6').Value = 163
$ws.Range('K7').Value = 684
$ws.Range('K8').Value = 1497
$ws.Range('K9').Value = 105
$ws.Range('K11').Value = 420
$ws.Range('K17').Value = 43
$ws.Range('K18').Value = 151
$ws.Range('K19').Value = 665
$ws.Range('K20').Value = 548
$ws.Range('K21').Value = 74
$ws.Range('K23').Value = 228
$ws.Range('K24').Value = 72
$ws.Range('K25').Value = 107
$ws.Range('K27').Value = 209
$ws.Range('K29').Value = 1238
$ws.Range('K31').Value = 252
$ws.Range('K37').Value = 775
$ws.Range('K42').Value = 842
$ws.Range('K44').Value = 189
$ws.Range('K46').Value = 48
$ws.Range('K48').Value = 288
$ws.Range('K52').Value = 606
$ws.Range('K53').Value = 289
$ws.Range('K54').Value = 451
$ws.Range('K55').Value = 246
$ws.Range('K59').Value = 41
$ws.Range('H63').Value = 292
$ws.Range('I63').Value = 228
$ws.Range('K63').Value = 60
$ws.Range('K64').Value = 142
$ws.Range('K67').Value = 891
$ws.Range('K73').Value = 207
$ws.Range('K76').Value = 306
$ws.Range('K77').Value = 154
$ws.Range('K78').Value = 258
$ws.Range('K83').Value = 494
$ws.Range('K85').Value = 1053
$ws.Range('K87').Value = 46
$ws.Range('K89').Value = 343
$ws.Range('K90').Value = 216
$ws.Range('K91').Value = 270
$ws.Range('K94').Value = 305
$ws.Range('K95').Value = 375
$ws.Range('K99').Value = 384
$ws.Range('H101').Value = 26053
$ws.Range('I101').Value = 26267
$ws.Range('K101').Value = 22827

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('K3').Value = 64
$ws.Range('K7').Value = 252

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('K6').Value = 253
$ws.Range('K7').Value = 891

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('K3').Value = 109
$ws.Range('K6').Value = 242
$ws.Range('K7').Value = 451

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('K4').Value = 59
$ws.Range('K6').Value = 363
$ws.Range('K7').Value = 1238

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('K6').Value = 136
$ws.Range('K7').Value = 288

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('K6').Value = 220
$ws.Range('K7').Value = 665

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range('K2').Value = 53
$ws.Range('K7').Value = 189

$ws = $wb.Worksheets.Item('River North')
$ws.Range('K2').Value = 70
$ws.Range('K7').Value = 306

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range('K2').Value = 61
$ws.Range('K7').Value = 163

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('K2').Value = 224
$ws.Range('K6').Value = 312
$ws.Range('K7').Value = 842

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('K5').Value = 7
$ws.Range('K7').Value = 258

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('K2').Value = 75
$ws.Range('K7').Value = 246

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range('K2').Value = 29
$ws.Range('K7').Value = 72

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range('K6').Value = 15
$ws.Range('K7').Value = 48

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range('K4').Value = 15
$ws.Range('K7').Value = 228

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('K2').Value = 69
$ws.Range('K7').Value = 270

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range('K6').Value = 45
$ws.Range('K7').Value = 74

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range('K6').Value = 54
$ws.Range('K7').Value = 142

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('K2').Value = 190
$ws.Range('K3').Value = 175
$ws.Range('K7').Value = 548

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range('K3').Value = 49
$ws.Range('K7').Value = 151

$ws = $wb.Worksheets.Item('Burnside')
$ws.Range('K2').Value = 17
$ws.Range('K7').Value = 43

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('K2').Value = 224
$ws.Range('K5').Value = 29
$ws.Range('K6').Value = 187
$ws.Range('K7').Value = 684

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('K6').Value = 138
$ws.Range('K7').Value = 305

$ws = $wb.Worksheets.Item('East Side')
$ws.Range('K3').Value = 36
$ws.Range('K7').Value = 107

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('K2').Value = 144
$ws.Range('K7').Value = 420

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range('K3').Value = 37
$ws.Range('K7').Value = 105

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('K2').Value = 71
$ws.Range('K6').Value = 70
$ws.Range('K7').Value = 207

$ws = $wb.Worksheets.Item('Montclare')
$ws.Range('K6').Value = 14
$ws.Range('K7').Value = 41

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('K3').Value = 106
$ws.Range('K7').Value = 343

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range('K3').Value = 50
$ws.Range('K7').Value = 209

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range('K2').Value = 79
$ws.Range('K7').Value = 216

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('K2').Value = 345
$ws.Range('K3').Value = 366
$ws.Range('K6').Value = 256
$ws.Range('K7').Value = 1053

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range('K2').Value = 64
$ws.Range('K7').Value = 154

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('K6').Value = 221
$ws.Range('K7').Value = 606

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Range('K6').Value = 22
$ws.Range('K7').Value = 46
